# Chip-8 opcode list: disassembler now supports all basic ch8 instructions.
# Adds a new "In Disassembler" column (D), shifting "In emulator" -> E and
# "Notes" -> F, marks every basic opcode with an "X" in the new column
# (except the ignored "0nnn SYS addr" opcode), and adds Notes for the
# SYS opcode plus the whole block of Super Chip-8 extension opcodes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: shift "In emulator" / "Notes" right, insert new column ---
$ws.Range("F1").Value = "Notes"                 # Notes -> F1
$ws.Range("E1").Value = "In emulator"            # In emulator -> E1
$ws.Range("D1").Value = "In Disassembler"

# --- Mark every basic Chip-8 opcode as supported in the disassembler ---
$disassemblerRows = 2,3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36
foreach ($r in $disassemblerRows) {
    $ws.Cells.Item($r, 4).Value = "X"
}

# --- Notes column (F) ---
$ws.Range("F4").Value = "Used by older computers and should be ignored"

$superChip8Rows = 37,38,39,40,41,42,43,44,45,46
foreach ($r in $superChip8Rows) {
    $ws.Cells.Item($r, 6).Value = "super chip8"
}

# --- Cosmetic: widen the new "In Disassembler" column, move the active cell ---
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Range("E10").Select() | Out-Null
